$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value of B2 from "super" to "supers"
$ws.Range("B2").Value = "supers"

# Move/set the active selection to B2 (matches the updated sheetView selection)
$ws.Range("B2").Select()
